$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" (D4) and "Correspond Handback DateTime" (G4)
# for the 0f1b92a9-af8d-4718-a057-8fe11696b002 entry on both the zh-cn and de-de sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-23 08:46:12"
$wsZhCn.Range("G4").Value = "2016-02-23 08:46:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-23 08:46:23"
$wsDeDe.Range("G4").Value = "2016-02-23 08:47:20"
